$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Shared text blocks reused across the new rows
$preDesc1 = 'convert to lower, convert unicode to ascii, remove multiple spaces, trim "space" and ","'
$modelDet1 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000'
$preDesc2 = 'remove multiple spaces, convert to lower, trim "space" and ",", convert unicode to ascii'
$modelDet2 = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000'

$features = '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit'
$model = 'Neuron Network'
$filter = '0 filters: '

$rows = @(
    @{ Row=23; A='20160405_164851'; B=2262.705; C=$preDesc1; D=$features; E=$model; F=$modelDet1; G=0.994;               H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=24; A='20160405_172633'; B=2328.383; C=$preDesc1; D=$features; E=$model; F=$modelDet1; G=0.993333333333333;   H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=25; A='20160405_180522'; B=2393.407; C=$preDesc1; D=$features; E=$model; F=$modelDet1; G=0.994666666666667;   H=0.99009900990099;   I=$filter; J=0.34375 },
    @{ Row=26; A='20160405_184515'; B=2424.495; C=$preDesc1; D=$features; E=$model; F=$modelDet1; G=0.992666666666667;   H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=27; A='20160405_192540'; B=2421.675; C=$preDesc1; D=$features; E=$model; F=$modelDet1; G=0.99;                H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=28; A='20160406_081417'; B=3407.352; C=$preDesc2; D=$features; E=$model; F=$modelDet2; G=0.994;               H=0.986798679867987;  I=$filter; J=0.368421052631579 },
    @{ Row=29; A='20160406_091105'; B=3422.952; C=$preDesc2; D=$features; E=$model; F=$modelDet2; G=0.991333333333333;   H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=30; A='20160406_100808'; B=5140.958; C=$preDesc2; D=$features; E=$model; F=$modelDet2; G=0.990666666666667;   H=0.99009900990099;   I=$filter; J=0.385416666666667 },
    @{ Row=31; A='20160406_113349'; B=5813.335; C=$preDesc2; D=$features; E=$model; F=$modelDet2; G=0.993333333333333;   H=0.99009900990099;   I=$filter; J=0.416666666666667 },
    @{ Row=32; A='20160406_131042'; B=7151.665; C=$preDesc2; D=$features; E=$model; F=$modelDet2; G=0.991333333333333;   H=0.99009900990099;   I=$filter; J=0.416666666666667 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
}
